# The roster had a duplicate entry for student 李世澤: row 6 used the
# (correct) student id "D1043618" but zero scores for 爬蟲1/爬蟲2/資料處裡_口頭報告,
# while row 7 used a mistyped lowercase id "d1043618" but carried the real
# scores for those same three columns. Consolidate the real scores onto the
# correctly-cased row, then remove the now-redundant duplicate row (the
# remaining rows shift up to close the gap).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("D6").Value = 100.0
$ws1.Range("E6").Value = 90.0
$ws1.Range("G6").Value = 100.0

$ws1.Rows.Item(7).Delete()

# The two other (previously blank) sheets each get the same header row as
# the main roster sheet, copied so the cell formatting (style) matches too.
$ws1.Range("A1:B1").Copy()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1").Value = "學號"
$ws2.Range("B1").Value = "姓名"
$ws2.Range("A1:B1").PasteSpecial(-4122)

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A1").Value = "學號"
$ws3.Range("B1").Value = "姓名"
$ws3.Range("A1:B1").PasteSpecial(-4122)
